# "Generate Report for Handback" - populate the handback columns (Latest
# Target File / Latest Handback File / Latest Handback DateTime) for the
# zh-cn and de-de localization targets, and flip the Overview/Status text
# from "Ready for handoff" to "Handed back: in sync with en-US".

$wb = $excel.ActiveWorkbook

$mdName    = "0b4c3046-7a2c-4471-a243-7a24bb35215f.md"
$mdUrl     = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/84bc305b75ba80199bf4a9ac6c369540a47ce1df/e2e/0b4c3046-7a2c-4471-a243-7a24bb35215f.md"
$zhTarget  = "0b4c3046-7a2c-4471-a243-7a24bb35215f.e37391fc87edcace882f353b43e03e8cfdd7e087.zh-cn.xlf"
$deTarget  = "0b4c3046-7a2c-4471-a243-7a24bb35215f.e37391fc87edcace882f353b43e03e8cfdd7e087.de-de.xlf"

# Every cell still showing "Ready for handoff" (Overview!E2/F2 and the
# Status column on the language sheets) is now reporting the synced state.
$ws1 = $wb.Worksheets.Item("Overview")
$ws1.Range("E2").Value = "Handed back: in sync with en-US"
$ws1.Range("F2").Value = "Handed back: in sync with en-US"

# zh-cn sheet: fill in target file / handback file / handback datetime.
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("C2").Value = "Handed back: in sync with en-US"
$wsZh.Range("I2").Value = $mdName
$wsZh.Hyperlinks.Add($wsZh.Range("I2"), $mdUrl, "", "", $mdName)
$wsZh.Range("J2").Value = $zhTarget
$wsZh.Range("K2").Value = "2016-08-17 06:52:28"

# de-de sheet: same, but its handback finished a few seconds later.
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("C2").Value = "Handed back: in sync with en-US"
$wsDe.Range("I2").Value = $mdName
$wsDe.Hyperlinks.Add($wsDe.Range("I2"), $mdUrl, "", "", $mdName)
$wsDe.Range("J2").Value = $deTarget
$wsDe.Range("K2").Value = "2016-08-17 06:52:35"

# The freshly-linked "Latest Target File" cells read like the other
# hyperlinked filename cells (underlined Calibri, matching the A2 style).
$wsZh.Range("I2").Font.Name = "Calibri"
$wsZh.Range("I2").Font.Size = 11
$wsZh.Range("I2").Font.Underline = 2
$wsZh.Range("I2").Font.Color = 15570276

$wsDe.Range("I2").Font.Name = "Calibri"
$wsDe.Range("I2").Font.Size = 11
$wsDe.Range("I2").Font.Underline = 2
$wsDe.Range("I2").Font.Color = 15570276

# Widen the columns that now carry the longer text so it's still readable.
$ws1.Columns.Item(5).ColumnWidth = 29.15
$ws1.Columns.Item(6).ColumnWidth = 29.15

$wsZh.Columns.Item(3).ColumnWidth = 29.15
$wsZh.Columns.Item(9).ColumnWidth = 39.15
$wsZh.Columns.Item(10).ColumnWidth = 39.15

$wsDe.Columns.Item(3).ColumnWidth = 29.15
$wsDe.Columns.Item(9).ColumnWidth = 39.15
$wsDe.Columns.Item(10).ColumnWidth = 39.15
